$wb = $excel.ActiveWorkbook

# Row 4 across all three sheets corresponds to file
# "640210b5-cd49-4c83-a319-599b5531f85b.md", whose status is
# "Ready for handoff". Generating the handoff report refreshes the
# handoff timestamps recorded for that file.

# Overview sheet: "Latest Handoff Date" column (D)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value = "2016-40-17 02:40:08"

# zh-cn sheet: "Latest Handoff Datetime" column (E)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-17 02:39:59"

# de-de sheet: "Latest Handoff Datetime" column (E)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-17 02:40:08"
